# Update the "取得日時" (retrieved at) timestamp for the newly appended
# case rows on the "ランサーズ" sheet from 12:39:26 to 12:52:26 (JST).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-23 12:52:26"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
